$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at position 4 for "A 18040-2021" (this record moves from
#    its old position - row 12 - up to row 4, with refreshed statistics), and
#    push all the following rows (old 4..307) down by one (new 5..308).
#    The record's old row (row 12, which after the insert above has shifted
#    down to row 13) is now a duplicate and must be removed.
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(13).Delete()

$ws.Range("A4").Value = "A 18040-2021"
$ws.Range("B4").Value = 44302
$ws.Range("C4").Value = 45192
$ws.Range("D4").Value = "DALARNAS LÄN"
$ws.Range("E4").Value = "HEDEMORA"
$ws.Range("G4").Value = 17.3
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 6
$ws.Range("P4").Value = 3
$ws.Range("Q4").Value = 13

$ws.Range("R4").Value = "Knärot`r`nRynkskinn`r`nUlltickeporing`r`nKortskaftad ärgspik`r`nTretåig hackspett`r`nUllticka`r`nBollvitmossa`r`nBronshjon`r`nMindre märgborre`r`nSkarp dropptaggsvamp`r`nTallfingersvamp`r`nVedticka`r`nVågbandad barkbock"

$ws.Range("S4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HEDEMORA/artfynd/A 18040-2021.xlsx", "A 18040-2021")'
$ws.Range("T4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HEDEMORA/kartor/A 18040-2021.png", "A 18040-2021")'
$ws.Range("U4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HEDEMORA/knärot/A 18040-2021.png", "A 18040-2021")'
$ws.Range("V4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HEDEMORA/klagomål/A 18040-2021.docx", "A 18040-2021")'
$ws.Range("W4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HEDEMORA/klagomålsmail/A 18040-2021.docx", "A 18040-2021")'
$ws.Range("X4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HEDEMORA/tillsyn/A 18040-2021.docx", "A 18040-2021")'
$ws.Range("Y4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HEDEMORA/tillsynsmail/A 18040-2021.docx", "A 18040-2021")'

# ---------------------------------------------------------------------------
# 2. Row 2 ("A 9053-2020") gets extra species found and its summary counts
#    updated.
# ---------------------------------------------------------------------------
$ws.Range("I2").Value = 11
$ws.Range("J2").Value = 10
$ws.Range("O2").Value = 13
$ws.Range("Q2").Value = 25

$ws.Range("R2").Value = "Knärot`r`nLammticka`r`nRynkskinn`r`nGarnlav`r`nGropig brunbagge`r`nGropticka`r`nGränsticka`r`nSpillkråka`r`nSvartvit flugsnappare`r`nTallticka`r`nTalltita`r`nUllticka`r`nViolettgrå tagellav`r`nBollvitmossa`r`nBronshjon`r`nGulnål`r`nRävticka`r`nSkinnlav`r`nStor aspticka`r`nSvavelriska`r`nThomsons trägnagare`r`nTibast`r`nVedticka`r`nVårärt`r`nVanlig groda"

# ---------------------------------------------------------------------------
# 3. Every report's "Förändrad" (changed) date moves from 2023-09-21 (45190)
#    to 2023-09-23 (45192). The net row count is unchanged (one row inserted,
#    one removed), so the data still spans rows 2..307.
# ---------------------------------------------------------------------------
$ws.Range("C2:C307").Value = 45192
